$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3 holds a date-like string ("2025-09-04"). Force text formatting first so
# Excel doesn't auto-convert it into a real date serial number, then clear
# the formatting override afterwards so the cell keeps the workbook's
# default (unstyled) look, matching a plain appended row.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-09-04"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = "15:19:52"
$ws.Range("C3").Value = "1.00 EUR = 1589.8516 ARS"
